$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45979
$ws.Range("B2").Value = 78.98999999999999
$ws.Range("C2").Value = 75.44
$ws.Range("D2").Value = 71.53
$ws.Range("E2").Value = 69.43000000000001
$ws.Range("F2").Value = 68.70999999999999
$ws.Range("G2").Value = 68.87
$ws.Range("H2").Value = 90.41
$ws.Range("I2").Value = 103.04
$ws.Range("J2").Value = 103.5
$ws.Range("K2").Value = 76.09
$ws.Range("L2").Value = 42.02
$ws.Range("M2").Value = 13.95
$ws.Range("N2").Value = 16.25
$ws.Range("O2").Value = 19.63
$ws.Range("P2").Value = 25.65
$ws.Range("Q2").Value = 48.83
$ws.Range("R2").Value = 82.73
$ws.Range("S2").Value = 106.06
$ws.Range("T2").Value = 131.68
$ws.Range("U2").Value = 133.05
$ws.Range("V2").Value = 134.38
$ws.Range("W2").Value = 134.42
$ws.Range("X2").Value = 112.65
$ws.Range("Y2").Value = 97.76000000000001
$ws.Range("Z2").Value = 79.38
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 119.8
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 134.4
$ws.Range("AE2").Value = "18h-20h"
$ws.Range("AF2").Value = 132.36
$ws.Range("AG2").Value = "0h-15h"
